$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: numeric summary row for "Empresa Modelo"
$ws.Range("A2").Value = "Empresa Modelo"
$ws.Range("B2").Value = 148622.773203
$ws.Range("C2").Value = -373.65
$ws.Range("D2").Value = 148249.123203
$ws.Range("E2").Value = 1

# Row 3: formatted-text mirror row
$ws.Range("B3").Value = "148.622,77"
$ws.Range("C3").Value = "-373,65"
$ws.Range("D3").Value = "148.249,12"
